# 11012017 - implementacao de alteracoes no formulario de configuracao,
# implementado verificacao de espaco no vinculo de tabela
#
# Populates row 8 of the monitoring table with the data-type legend used
# to validate/annotate each column (Codigo / float/10 / int / float /
# "v1.v2.v3.v4.v5"), matching the formatting already used by the
# neighbouring header cells (text number format, centered).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlCenter = -4108

# A8 - row label
$ws.Range("A8").Value = "Codigo"

# B8:G8 and I8:N8 - "float/10" columns (IMEI/SIM info fields)
$floatTenCells = @("B8","C8","D8","E8","F8","G8","I8","J8","K8","L8","M8","N8")
foreach ($addr in $floatTenCells) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.HorizontalAlignment = $xlCenter
    $c.VerticalAlignment = $xlCenter
    $c.Value = "float/10"
}

# H8 and O8 - "int" (merged with row 7 above, SIM field)
$intCells = @("H8","O8")
foreach ($addr in $intCells) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.HorizontalAlignment = $xlCenter
    $c.VerticalAlignment = $xlCenter
    $c.Value = "int"
}

# P8 - "float"
$p8 = $ws.Range("P8")
$p8.NumberFormat = "@"
$p8.HorizontalAlignment = $xlCenter
$p8.VerticalAlignment = $xlCenter
$p8.Value = "float"

# Q8 - "float " (trailing space preserved)
$q8 = $ws.Range("Q8")
$q8.NumberFormat = "@"
$q8.HorizontalAlignment = $xlCenter
$q8.VerticalAlignment = $xlCenter
$q8.Value = "float "

# R8:T8 - version format
$versionCells = @("R8","S8","T8")
foreach ($addr in $versionCells) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.HorizontalAlignment = $xlCenter
    $c.VerticalAlignment = $xlCenter
    $c.Value = "v1.v2.v3.v4.v5"
}

# Widen columns L:N and R:T to fit the newly typed values.
$ws.Columns("L:N").ColumnWidth = 7.1
$ws.Columns("R:T").ColumnWidth = 12.6

# Restore the cursor/selection to the merged cell O4:O5, as left by the author.
$null = $ws.Range("O4:O5").Select()
